# STM32 Pin Allocation.xlsx edit
# - Re-optimized pin allocation: several "Feature/Notes" labels in columns G and P
#   were reassigned between rows (TVS diode / current-sense / relay / ADC
#   labelling reshuffle).
# - Row 37 grew taller to fit the now-longer "Output Current Sense 1 (ADC)" text.
# - Selection moved to where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (left bank) relabeling
$ws.Range("G2").Value  = "I_OUT_1_OCD (GPIO)"
$ws.Range("G3").Value  = "I_IN_1_OCD (GPIO)"
$ws.Range("G4").Value  = "I_OUT_2_OCD (GPIO)"
$ws.Range("G5").Value  = "I_IN_2_OCD (GPIO)"
$ws.Range("G22").Value = "I_OUT_3_OCD (GPIO)"
$ws.Range("G36").Value = "Output Current Sense 1 (ADC)"
$ws.Range("G37").Value = "Input Current Sense 3 (ADC)"

# Column P (right bank) relabeling
$ws.Range("P18").Value = "Ext_ADC_EOC (GPIO)"
$ws.Range("P19").Value = "Relay Control 2 (GPIO)"
$ws.Range("P20").Value = "I_IN_3_OCD (GPIO)"
$ws.Range("P23").Value = "Relay Control 3 (GPIO)"

# Row 37 now needs extra height for the longer wrapped label
$ws.Rows.Item(37).RowHeight = 43.15

# Leave the selection where the author ended up working
$ws.Range("P30").Select()
